# Update cryptos list values (price + 1h volume change) per latest scrape.
# Rows 46/47 swap rank: InjectiveProtocol now ranks above Kaspa.
# Price cells whose text would otherwise be auto-parsed as a number are
# forced to Text format first so values like trailing-zero prices
# ("5.80", "10.90") and exact decimals are preserved verbatim as strings,
# matching the source data (Price column is textual, not numeric).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.030.92"
$ws.Range("E2").Value = "  -0.36%  "
$ws.Range("D3").Value = "1.778.80"
$ws.Range("E3").Value = "  -2.58%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "225.61"
$ws.Range("E5").Value = "  +0.28%  "
$ws.Range("E6").Value = "  -1.53%  "
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.41"
$ws.Range("E8").Value = "  +1.41%  "
$ws.Range("E9").Value = "  -2.16%  "
$ws.Range("E10").Value = "  -1.94%  "
$ws.Range("E11").Value = "  +0.53%  "
$ws.Range("D12").Value = "2.038.00"
$ws.Range("E12").Value = "  -2.51%  "
$ws.Range("D13").Value = "1.866.13"
$ws.Range("E13").Value = "  +2.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.90"
$ws.Range("E14").Value = "  +0.56%  "
$ws.Range("D15").Value = "33.960.25"
$ws.Range("E15").Value = "  -0.56%  "
$ws.Range("E16").Value = "  -4.21%  "
$ws.Range("E17").Value = "  -5.06%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "67.64"
$ws.Range("E18").Value = "  -2.87%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "243.33"
$ws.Range("E19").Value = "  -3.12%  "
$ws.Range("E21").Value = "  +0.08%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.67"
$ws.Range("E22").Value = "  -4.23%  "
$ws.Range("E23").Value = "  -4.90%  "
$ws.Range("E24").Value = "  -4.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "159.95"
$ws.Range("E25").Value = "  -0.44%  "
$ws.Range("E26").Value = "  -2.36%  "
$ws.Range("E27").Value = "  -3.15%  "
$ws.Range("E28").Value = "  -2.41%  "
$ws.Range("E29").Value = "  +0.16%  "
$ws.Range("E30").Value = "  +0.65%  "
$ws.Range("E31").Value = "  -4.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.61"
$ws.Range("E32").Value = "  -4.40%  "
$ws.Range("E33").Value = "  -2.51%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.79"
$ws.Range("E34").Value = "  -5.27%  "
$ws.Range("D35").Value = "1.386.11"
$ws.Range("E35").Value = "  -3.64%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.645"
$ws.Range("E36").Value = "  +0.49%  "
$ws.Range("E37").Value = "  -2.06%  "
$ws.Range("E38").Value = "  -2.09%  "
$ws.Range("E39").Value = "  +0.19%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.19"
$ws.Range("E40").Value = "  +1.57%  "
$ws.Range("E41").Value = "  -2.92%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.908"
$ws.Range("E42").Value = "  -5.56%  "
$ws.Range("E43").Value = "  -5.12%  "
$ws.Range("E44").Value = "  +13.35%  "
$ws.Range("E45").Value = "  +1.78%  "
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.46"
$ws.Range("E46").Value = "  +4.70%  "
$ws.Range("B47").Value = "Kaspa"
$ws.Range("C47").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0497"
$ws.Range("E47").Value = "  -0.17%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "107.39"
$ws.Range("E48").Value = "  +0.47%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.80"
$ws.Range("E49").Value = "  -4.66%  "
$ws.Range("D50").Value = "1.936.58"
$ws.Range("E50").Value = "  -2.52%  "
$ws.Range("E51").Value = "  +0.11%  "
